$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (shared text) on both sheets
$newProductName = "4210-RBI-EI-DB-DL-REC-RNI-INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ONT-PER-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update shortname on input sheet: was numeric 4210, now text "421q"
$wsInput.Range("B2").Value = "421q"

# Reset the input sheet's view: scroll back to top and select B1
$wsInput.Activate()
$wsInput.Application.ActiveWindow.ScrollRow = 1
$wsInput.Application.ActiveWindow.ScrollColumn = 1
$wsInput.Range("B1").Select()
